$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: sexMale
$ws.Range("A2").Value = "sexMale"
$ws.Range("B2").Value = 1.085687007161959
$ws.Range("C2").Value = 0.9164260556025183
$ws.Range("D2").Value = 1.286209913297725
$ws.Range("E2").Value = 0.3417468030667726

# Row 3: age
$ws.Range("A3").Value = "age"
$ws.Range("B3").Value = 1.034481121476051
$ws.Range("C3").Value = 1.027853885834152
$ws.Range("D3").Value = 1.041151087172152
$ws.Range("E3").Value = [double]"4.73606371603818e-25"
